# The "Recorded By" column (G) on the "Session Analysis Results" sheet lists
# the users/processes that recorded or touched each attendance session, e.g.
# "System, dnasr281@gmail.com" or "System, system, backup@backdoor.com".
#
# This sync moves the literal, capitalised "System" entry from the front of
# that comma-separated list to the end, e.g.:
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com"     -> "system, backup@backdoor.com, System"
#
# Rows whose value does not start with "System, " (blank cells, cells that
# are just "System" on their own, or cells that already have "System"
# somewhere other than first) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2
    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text.StartsWith("System, ")) {
        $parts = New-Object System.Collections.Generic.List[string]
        foreach ($p in $text.Split(",")) {
            [void]$parts.Add($p.Trim())
        }
        if ($parts[0] -eq "System") {
            [void]$parts.RemoveAt(0)
            [void]$parts.Add("System")
            $cell.Value = [string]::Join(", ", $parts)
        }
    }
}
